$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add AUC (Area under ROC) values for "Com o dataset da Feature Reduction" section (row 13)
$ws.Range("B13").Value = 0.881
$ws.Range("C13").Value = 0.957
$ws.Range("D13").Value = 0.956

# Add AUC (Area under ROC) values for "Com o dataset da Feature Selection" section (row 23)
$ws.Range("B23").Value = 0.831
$ws.Range("C23").Value = 0.761
$ws.Range("D23").Value = 0.761

# Update the view state (scroll position / selection) to match the saved workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select()
